# Automatische test-sync: 2025-08-28 21:02:50
# Add a new "Retour status" log row to the Logs sheet, extend the
# conditional-formatting ranges to cover it, and bump the
# "Retour / Terugbetaling" count on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$oldLastRow = 18
$newLastRow = 19

# --- Append the new row of log data -----------------------------------
$logs.Cells.Item($newLastRow, 1).Value = "Retour status"
$logs.Cells.Item($newLastRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newLastRow, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($newLastRow, 6).Value = "2025-08-28 21:02:06"
$logs.Cells.Item($newLastRow, 7).Value = "Ja"
$logs.Cells.Item($newLastRow, 8).Value = "Nee"
$logs.Cells.Item($newLastRow, 9).Value = "Nee"
$logs.Cells.Item($newLastRow, 10).Value = "Nee"

# --- Extend the conditional formatting sqref ranges by one row --------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + $oldLastRow)
    $newRange = $logs.Range($col + "2:" + $col + $newLastRow)
    $cfs = $oldRange.FormatConditions
    for ($i = 1; $i -le $cfs.Count; $i++) {
        $cfs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard summary count --------------------------------
$dashboard.Range("B2").Value = 17
